# "subimos el último SPA" - update last week's target units (column R) for the
# rows affected by the new SPA (Semana Pasada) figures, then recompute the
# dependent columns:
#   T (Tendencia Consumo, 20=19-18)   = MAX(0, S - R)
#   U (Pedido Final, 21=6+11-16+20)   = MAX(0, F + K - P + T)
# and finally refresh the "Total_Unidades" summary metric in C43, which is
# the sum of the "Pedido Final" column (U) over the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "uds. Objetivo semana pasada" (column R) values, keyed by row number.
$rowsToUpdate = @{
    5  = 2
    6  = 1
    10 = 1
    11 = 3
    12 = 2
    15 = 1
    16 = 1
    17 = 3
    23 = 4
    26 = 1
    27 = 2
    34 = 4
    36 = 1
    37 = 3
    38 = 7
    39 = 4
}

foreach ($row in $rowsToUpdate.Keys) {
    $ws.Cells.Item($row, 18).Value2 = $rowsToUpdate[$row]
}

# Data rows run from 3 to 40. Recompute the dependent columns for every row
# so the sheet stays internally consistent (rows whose inputs didn't change
# recompute to the same value they already had).
# NOTE: plain `.Value` reads are unreliable in this COM shim, so `.Value2`
# is used for every get/set below.
$firstDataRow = 3
$lastDataRow = 40
$totalUnidades = 0

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $fVal = $ws.Cells.Item($row, 6).Value2   # F: Unidades Calculadas
    $kVal = $ws.Cells.Item($row, 11).Value2  # K: Stock Minimo Objetivo
    $pVal = $ws.Cells.Item($row, 16).Value2  # P: Stock Real
    $rVal = $ws.Cells.Item($row, 18).Value2  # R: uds. Objetivo semana pasada
    $sVal = $ws.Cells.Item($row, 19).Value2  # S: Uds. Vtas. reales semana pasada

    if ($null -eq $fVal) { $fVal = 0 }
    if ($null -eq $kVal) { $kVal = 0 }
    if ($null -eq $pVal) { $pVal = 0 }
    if ($null -eq $rVal) { $rVal = 0 }
    if ($null -eq $sVal) { $sVal = 0 }

    $tVal = $sVal - $rVal
    if ($tVal -lt 0) { $tVal = 0 }

    $uVal = $fVal + $kVal - $pVal + $tVal
    if ($uVal -lt 0) { $uVal = 0 }

    $ws.Cells.Item($row, 20).Value2 = $tVal   # T: Tendencia Consumo
    $ws.Cells.Item($row, 21).Value2 = $uVal   # U: Pedido Final

    $totalUnidades = $totalUnidades + $uVal
}

# C43 = Total_Unidades, the sum of column U (Pedido Final) across data rows.
$ws.Cells.Item(43, 3).Value2 = $totalUnidades
